$d = $word.ActiveDocument

# Locate the unique sentence-ending run that reads:
#   "...Add this to your main.ts file:"
# and grab just the " file:" tail (the two words after the highlighted
# "main.ts" run, up to and including the trailing colon).
$target = $d.Content
$found = $target.Find.Execute("main.ts file:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'main.ts file:' text"
}

$fileColon = $d.Range($target.End - 6, $target.End)
if ($fileColon.Text -ne " file:") {
    throw "Unexpected text in target range: [$($fileColon.Text)]"
}

# " file:" -> " file" (drop the trailing colon; it gets re-added at the
# very end of the new sentence tail below).
$fileColon.Text = " file"

# Append the rest of the new sentence as separate runs so "shell" can
# carry its own (highlighted) run properties, matching how the other
# inline code references (main.ts, bootstrap.ts, loadRemoteEntry) are
# marked up elsewhere in this document.
$r1 = $d.Range($fileColon.End, $fileColon.End)
$r1.InsertAfter(" in your ")

$r2 = $d.Range($r1.End, $r1.End)
$r2.InsertAfter("shell")
$r2.Font.HighlightColorIndex = 15   # wdGray50 -> <w:highlight w:val="darkGray"/>

$r3 = $d.Range($r2.End, $r2.End)
$r3.InsertAfter(" project")

$r4 = $d.Range($r3.End, $r3.End)
$r4.InsertAfter(":")

Write-Output "Updated 'main.ts file:' sentence to mention the shell project"
